# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Murcott, Primera & Segunda, 2022-12-23)
# right before the current row 157, shifting the existing rows 157:180 down
# to 159:182 (dimension grows from A1:T180 to A1:T182).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 157 (pushes old 157.. down by 2)
$ws.Rows("157:158").Insert()

# --- New row 157: Murcott / Primera ---
$ws.Range("A157").Value = 11
$ws.Range("B157").Value = "Vega Monumental Concepción"
$ws.Range("C157").Value = "Bíobío"
$ws.Range("D157").Value = 44918
$ws.Range("E157").Value = 8
$ws.Range("F157").Value = "Fruta"
$ws.Range("G157").Value = 100102
$ws.Range("H157").Value = "Cítricos"
$ws.Range("I157").Value = 100102004
$ws.Range("J157").Value = "Mandarina"
$ws.Range("K157").Value = "Murcott"
$ws.Range("L157").Value = "Primera"
$ws.Range("M157").Value = 100
$ws.Range("N157").Value = 9000
$ws.Range("O157").Value = 10000
$ws.Range("P157").Value = 9500
$ws.Range("Q157").Value = "$/bandeja 18 kilos"
$ws.Range("R157").Value = "Región de O'Higgins"
$ws.Range("S157").Value = 528
$ws.Range("T157").Value = 18

# --- New row 158: Murcott / Segunda ---
$ws.Range("A158").Value = 11
$ws.Range("B158").Value = "Vega Monumental Concepción"
$ws.Range("C158").Value = "Bíobío"
$ws.Range("D158").Value = 44918
$ws.Range("E158").Value = 8
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100102
$ws.Range("H158").Value = "Cítricos"
$ws.Range("I158").Value = 100102004
$ws.Range("J158").Value = "Mandarina"
$ws.Range("K158").Value = "Murcott"
$ws.Range("L158").Value = "Segunda"
$ws.Range("M158").Value = 50
$ws.Range("N158").Value = 8000
$ws.Range("O158").Value = 8000
$ws.Range("P158").Value = 8000
$ws.Range("Q158").Value = "$/bandeja 18 kilos"
$ws.Range("R158").Value = "Región de O'Higgins"
$ws.Range("S158").Value = 444
$ws.Range("T158").Value = 18
